$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update weight values (column D) for rows 3, 4, 5, 7, 8
$ws.Range("D3").Value = 0.27
$ws.Range("D4").Value = 0.28
$ws.Range("D5").Value = 0.27
$ws.Range("D7").Value = 0.24
$ws.Range("D8").Value = 0.21

# Update the irl-price formula multiplier from 0.025 to 0.027 in column S
# (S3 has its own formula; S4:S8 share one formula group) - set the whole
# range at once so every row keeps its own relative references.
$ws.Range("S3").Formula = "=ROUND(Q3*0.027+P3+R3, 2)"
$ws.Range("S4:S8").Formula = "=ROUND(Q4*0.027+P4+R4, 2)"

# Move the active selection from F7 to D7
$ws.Range("D7").Select()
